# Add a new "Sheet2" worksheet between "Sheet1" and "Sheet1 (2)" that
# reports the simtime (s) used for each model's runtime benchmark.
#
# We build the new sheet by copying Sheet1 (so it inherits the same
# row-height / page-margin formatting) immediately after Sheet1, then
# clear its contents and repopulate it with the new table.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)

# Duplicate Sheet1 and drop the copy right after Sheet1; this becomes
# our new "Sheet2".
$sheet1.Copy($null, $sheet1)
$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Name = "Sheet2"
$sheet2.Cells.Clear()

# Header row
$sheet2.Range("A1").Value = "simtime (s)"
$sheet2.Range("B1").Value = "model"
$sheet2.Range("C1").Value = "full (s)"
$sheet2.Range("D1").Value = "particle (s)"
$sheet2.Range("E1").Value = "pct_speedup (%)"

# FK row
$sheet2.Range("A2").Value = 1
$sheet2.Range("B2").Value = "FK"
$sheet2.Range("C2").Value = 127
$sheet2.Range("D2").Value = 0.986
$sheet2.Range("E2").Formula = "=(C2-D2)/D2*100"

# LR row
$sheet2.Range("A3").Value = 0.1
$sheet2.Range("B3").Value = "LR"
$sheet2.Range("C3").Value = 31.4
$sheet2.Range("D3").Value = 0.119
$sheet2.Range("E3").Formula = "=(C3-D3)/D3*100"

# Update selections / active states to match the new layout:
#  - Sheet1 is no longer the active tab; its selection becomes the
#    whole data range with no specific active cell.
#  - Sheet2 becomes the active tab, with B1 selected.
$sheet1.Range("A1:D3").Select()
$sheet2.Range("B1").Select()
